# RPA datasets push 2024-05-18
# Insert two new IPO-underwriting records into the Sheet1 table:
#   - KB제28호스팩 (KB), subscription 2024-05-07 / payment 2024-05-10 / listing 2024-05-17
#   - 아이씨티케이 (NH), subscription 2024-05-07 / payment 2024-05-10 / listing 2024-05-17
# All pre-existing rows keep their values and simply shift down to make room.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force the cell to hold a literal text value (not an auto-converted
    # date/number) while leaving it with the default ("Normal") style, same
    # as every other data cell in the sheet.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Set-NumberCell($row, $col, $num) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = $num
    $cell.Style = "Normal"
}

# --- Insert the new KB row right after the header row (becomes row 2) ---
$ws.Rows.Item(2).Insert()

Set-TextCell 2 1 "KB"
Set-TextCell 2 2 "2024-05-07"
Set-TextCell 2 3 "KB제28호스팩"
Set-TextCell 2 4 "KB"
Set-TextCell 2 5 "KB"
Set-TextCell 2 6 "2024-05-10"
Set-TextCell 2 7 "2024-05-17"
Set-NumberCell 2 8 10000
Set-NumberCell 2 9 5000000
Set-NumberCell 2 10 2000
Set-NumberCell 2 11 0
Set-NumberCell 2 12 100

# --- Insert the new NH row right after the existing NH / 엔젤로보틱스 row ---
# Before this second insert, that row sits at row 6 (header=1, new KB row=2,
# old rows 2-4 now at 3-5, old row 5 "엔젤로보틱스" now at row 6).
$ws.Rows.Item(7).Insert()

Set-TextCell 7 1 "NH"
Set-TextCell 7 2 "2024-05-07"
Set-TextCell 7 3 "아이씨티케이"
Set-TextCell 7 4 "NH"
Set-TextCell 7 5 "NH"
Set-TextCell 7 6 "2024-05-10"
Set-TextCell 7 7 "2024-05-17"
Set-NumberCell 7 8 39400
Set-NumberCell 7 9 1970000
Set-NumberCell 7 10 20000
Set-NumberCell 7 11 0
Set-NumberCell 7 12 100
